$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C1 and C2 values
$ws.Range("C1").Value = 0.02552967818679151
$ws.Range("C2").Value = 0.03814844487578101

# Add new D and E columns
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 212.3515955309401

$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 199.8752016571346
